$d = $word.ActiveDocument

$oldText = "The data you’ll be using focuses on players that were selected in the first round of the NBA draft between the years 1990-2021, and they are divided based on what number in the first round they were selected"

$range = $d.Content
$range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$start = $range.Start

# Segments of the replacement text
$seg1 = "The data you’ll be using comes from the dataset "
$seg2 = "nba_draft.csv"
$seg3 = " "
$seg4 = "and "
$seg5 = "focuses on players that were selected in the first round of the NBA draft between the years 1990-2021, and they are divided based on what number in the first round they were selected"

$newText = $seg1 + $seg2 + $seg3 + $seg4 + $seg5

# Replace the whole run's text in one shot, keeping the original (minorHAnsi) formatting
$range.Text = $newText

# Compute character offsets (relative to $start) for each segment
$off1 = 0
$off2 = $off1 + $seg1.Length
$off3 = $off2 + $seg2.Length
$off4 = $off3 + $seg3.Length
$off5 = $off4 + $seg4.Length
$offEnd = $off5 + $seg5.Length

# "nba_draft.csv" -> Courier New, bold
$r2 = $d.Range($start + $off2, $start + $off3)
$r2.Font.Name = "Courier New"
$r2.Font.Bold = $true

# " " (space after the filename) -> Courier New, not bold
$r3 = $d.Range($start + $off3, $start + $off4)
$r3.Font.Name = "Courier New"
$r3.Font.Bold = $false
